$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new parameter rows -------------------------------------------------
# 1) New row 11: A_B_spacing (pushes A_joint_tolerance and everything below down by one)
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = "A_B_spacing"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "mm"

# 2) New row 15: B_height (right after B_length, pushes B_width.. down by one)
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "B_height"
$ws.Range("B15").Value = 50
$ws.Range("C15").Value = "mm"

# 3) New row 20: B_hinge_holes_spacing (right after B_hole_height, pushes B_hole_length.. down by one)
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "B_hinge_holes_spacing"
$ws.Range("B20").Value = 50
$ws.Range("C20").Value = "mm"

# --- Update existing values -----------------------------------------------------
# B_length value change 1200 -> 700 (now at row 14)
$ws.Range("B14").Value = 700

# C_length value change 300 -> 400 (now at row 23)
$ws.Range("B23").Value = 400

# C_height value change 400 -> 330 (now at row 25)
$ws.Range("B25").Value = 330

# --- Restore the view/selection state ------------------------------------------
[void]$ws.Range("B14").Select()
